$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

# Row 2
Set-TextValue 2 4 '34.216.19'
Set-TextValue 2 5 '  -1.08%  '

# Row 3
Set-TextValue 3 4 '1.783.96'
Set-TextValue 3 5 '  -2.39%  '

# Row 4
Set-TextValue 4 5 '  +0.22%  '

# Row 5
Set-TextValue 5 4 '224.82'
Set-TextValue 5 5 '  -2.58%  '

# Row 6
Set-TextValue 6 5 '  +0.15%  '

# Row 7
Set-TextValue 7 5 '  +0.12%  '

# Row 8
Set-TextValue 8 4 '31.90'
Set-TextValue 8 5 '  +0.26%  '

# Row 9
Set-TextValue 9 5 '  -1.30%  '

# Row 10
Set-TextValue 10 5 '  -2.24%  '

# Row 11
Set-TextValue 11 5 '  -0.15%  '

# Row 12
Set-TextValue 12 4 '2.035.32'
Set-TextValue 12 5 '  -2.50%  '

# Row 13
Set-TextValue 13 4 '11.19'
Set-TextValue 13 5 '  +7.20%  '

# Row 14
Set-TextValue 14 4 '1.787.16'
Set-TextValue 14 5 '  -2.21%  '

# Row 15
Set-TextValue 15 5 '  -4.05%  '

# Row 16
Set-TextValue 16 4 '34.211.51'
Set-TextValue 16 5 '  -0.69%  '

# Row 17
Set-TextValue 17 4 '4.22'
Set-TextValue 17 5 '  -1.86%  '

# Row 18
Set-TextValue 18 4 '68.80'
Set-TextValue 18 5 '  -1.37%  '

# Row 19
Set-TextValue 19 4 '254.75'
Set-TextValue 19 5 '  -1.77%  '

# Row 20
Set-TextValue 20 4 '0.0₃0740'
Set-TextValue 20 5 '  -1.95%  '

# Row 21
Set-TextValue 21 5 '  +0.05%  '

# Row 22
Set-TextValue 22 4 '10.37'
Set-TextValue 22 5 '  -2.01%  '

# Row 23
Set-TextValue 23 4 '4.21'
Set-TextValue 23 5 '  -3.67%  '

# Row 24
Set-TextValue 24 4 '2.13'
Set-TextValue 24 5 '  -4.34%  '

# Row 25
Set-TextValue 25 4 '157.56'
Set-TextValue 25 5 '  -0.89%  '

# Row 26
Set-TextValue 26 4 '16.41'
Set-TextValue 26 5 '  -1.80%  '

# Row 27
Set-TextValue 27 4 '7.02'
Set-TextValue 27 5 '  -1.77%  '

# Row 28
Set-TextValue 28 5 '  -1.28%  '

# Row 29
Set-TextValue 29 4 '0.999'
Set-TextValue 29 5 '  -0.16%  '

# Row 30
Set-TextValue 30 4 '3.78'
Set-TextValue 30 5 '  -2.99%  '

# Row 31
Set-TextValue 31 4 '0.0514'
Set-TextValue 31 5 '  -1.47%  '

# Row 32
Set-TextValue 32 5 '  -2.48%  '

# Row 33
Set-TextValue 33 4 '3.60'
Set-TextValue 33 5 '  +0.24%  '

# Row 34
Set-TextValue 34 4 '1.87'
Set-TextValue 34 5 '  +4.23%  '

# Row 35
Set-TextValue 35 4 '1.440.81'
Set-TextValue 35 5 '  -7.27%  '

# Row 36
Set-TextValue 36 5 '  -3.38%  '

# Row 37
Set-TextValue 37 4 '0.0188'
Set-TextValue 37 5 '  -1.40%  '

# Row 38
Set-TextValue 38 4 '0.625'
Set-TextValue 38 5 '  -1.92%  '

# Row 39
Set-TextValue 39 2 'MXToken'
Set-TextValue 39 3 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 39 4 '2.86'
Set-TextValue 39 5 '  +1.71%  '

# Row 40
Set-TextValue 40 2 'Aave'
Set-TextValue 40 3 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 40 4 '83.04'
Set-TextValue 40 5 '  -2.26%  '

# Row 41
Set-TextValue 41 5 '  +0.92%  '

# Row 42
Set-TextValue 42 4 '0.891'
Set-TextValue 42 5 '  -3.16%  '

# Row 43
Set-TextValue 43 5 '  -5.42%  '

# Row 44
Set-TextValue 44 5 '  -2.67%  '

# Row 45
Set-TextValue 45 5 '  -1.94%  '

# Row 46
Set-TextValue 46 2 'RocketPoolETH'
Set-TextValue 46 3 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue 46 4 '1.939.04'
Set-TextValue 46 5 '  -2.40%  '

# Row 47
Set-TextValue 47 2 'FraxShare'
Set-TextValue 47 3 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 47 4 '5.82'
Set-TextValue 47 5 '  +0.42%  '

# Row 48
Set-TextValue 48 4 '12.24'
Set-TextValue 48 5 '  -2.75%  '

# Row 49
Set-TextValue 49 4 '0.998'
Set-TextValue 49 5 '  +0.00%  '

# Row 50
Set-TextValue 50 4 '98.42'
Set-TextValue 50 5 '  +0.52%  '

# Row 51
Set-TextValue 51 4 '49.51'
Set-TextValue 51 5 '  -7.00%  '
